# devpo 19460 hotfix: refresh battery datapoint abbreviation table
# (B_CHARGE_LEVEL -> full B_* abbreviation catalogue, rows 11:56, col B widened)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Field Type (col A), Name Attribute (col B), Unit (col C), Description (col D)
# $null means "leave/clear that cell" (no unit / no description for that datapoint).
$rows = @(
    @(11, "datapoints", "B_CAPACITY", "Ah", "Nominal capacity"),
    @(12, "datapoints", "B_E_CHARGE_AC", "kWh", "Chargeable Energy"),
    @(13, "datapoints", "B_E_DISCHARGE_AC", "kWh", "Dischargeable Energy"),
    @(14, "datapoints", "B_E_EXP", "kWh", "Total energy from stacks"),
    @(15, "datapoints", "B_E_EXP_AC", "kWh", "Energy export from storage system AC"),
    @(16, "datapoints", "B_E_IMP", "kWh", "Total energy to stacks"),
    @(17, "datapoints", "B_E_IMP_AC", "kWh", "Energy import to storage system AC"),
    @(18, "datapoints", "B_E_INT_EXP", "kWh", "erzeugte Energie pro Intervall (geliefert)"),
    @(19, "datapoints", "B_E_INT_IMP", "kWh", "erzeugte Energie pro Intervall (bezogen)"),
    @(20, "datapoints", "B_E_STORED", "kWh", "Currently stored energy"),
    @(21, "datapoints", "B_F_AC", "Hz", "Grid frequency"),
    @(22, "datapoints", "B_I_AC", "A", "Battery AC current"),
    @(23, "datapoints", "B_I_DC", "A", "Ladestrom (DC)"),
    @(24, "datapoints", "B_LIM_I_CHARGE", $null, $null),
    @(25, "datapoints", "B_LIM_I_DISCHARGE", $null, $null),
    @(26, "datapoints", "B_LIM_P_CHARGE", "W", "Maximum charging power"),
    @(27, "datapoints", "B_LIM_P_DISCHARGE", "W", "Maximum discharging power"),
    @(28, "datapoints", "B_LIM_U_CHARGE", "V", "Charge end voltage"),
    @(29, "datapoints", "B_LIM_U_DISCHARGE", $null, $null),
    @(30, "datapoints", "B_OT_TOTAL", "h", "Operating Hours"),
    @(31, "datapoints", "B_P_AC", "W", "Battery power AC"),
    @(32, "datapoints", "B_P_DC", "W", "Total battery power"),
    @(33, "datapoints", "B_Q_AC", "VAr", "Battery reactive power AC"),
    @(34, "datapoints", "B_SOC", "%", "State of Charge"),
    @(35, "datapoints", "B_SOCH", "%", "Ladezustand (Nennkapazität) in %"),
    @(36, "datapoints", "B_SOH", "%", "Alterungsbedingter Erhaltungszustand in %"),
    @(37, "datapoints", "B_S_AC", "VA", "Battery apparent power AC"),
    @(38, "datapoints", "B_T_CELL[1..x]_[1..x]_[1..x]", "°C", "Cell temperature  [°C]"),
    @(39, "datapoints", "B_T_CELL_MAX[1..x]_[1..x]_[1..x]", "°C", "Maximum cell temperature  [°C]"),
    @(40, "datapoints", "B_T_CELL_MIN[1..x]_[1..x]_[1..x]", "°C", "Minimum cell temperature  [°C]"),
    @(41, "datapoints", "B_T_M[1..x]", "°C", "Module temperature [°C]"),
    @(42, "datapoints", "B_T_M_MAX[1..x]", "°C", "Maximum module temperature [°C]"),
    @(43, "datapoints", "B_T_M_MIN[1..x]", "°C", "Minimum module temperature [°C]"),
    @(44, "datapoints", "B_T_U[1..x]", "°C", "Temperature Outside/ Ambient [°C]"),
    @(45, "datapoints", "B_U_AC", "V", "Battery AC voltage"),
    @(46, "datapoints", "B_U_BULK", "V", "Battery charging voltage DC"),
    @(47, "datapoints", "B_U_CELL_AVG", "V", "Average cell voltage"),
    @(48, "datapoints", "B_U_CELL_MAX[1..x]_[1..x]_[1..x]", "V", "Maximum cell voltage"),
    @(49, "datapoints", "B_U_CELL_MIN[1..x]_[1..x]_[1..x]", "V", "Minimum cell voltage"),
    @(50, "datapoints", "B_U_DC", "V", "Battery voltage"),
    @(51, "datapoints", "B_U_OC", "V", "Open circuit voltage"),
    @(52, "datapoints", "T[1..x]", "°C", "Temperatures"),
    @(53, "datapoints", "STATE[1..x]", $null, "Global battery state conditions"),
    @(54, "datapoints", "ERROR[1..x]", $null, "Global battery error conditions"),
    @(55, "datapoints", "QS_TX", $null, "Telegrams transmitted (communication quality)"),
    @(56, "datapoints", "QS_RX", $null, "Telegrams received (communication quality)")
)

foreach ($entry in $rows) {
    $r = $entry[0]

    $ws.Range("A$r").Value = $entry[1]
    $ws.Range("B$r").Value = $entry[2]

    if ($null -eq $entry[3]) {
        $ws.Range("C$r").ClearContents() | Out-Null
    } else {
        $ws.Range("C$r").Value = $entry[3]
    }

    if ($null -eq $entry[4]) {
        $ws.Range("D$r").ClearContents() | Out-Null
    } else {
        $ws.Range("D$r").Value = $entry[4]
    }
}

# Column B now holds longer abbreviations (e.g. B_T_CELL_MAX[1..x]_[1..x]_[1..x]) - widen to fit.
$ws.Columns.Item(2).ColumnWidth = 37.92
